$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.975.45"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "2.098.77"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.81%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5156"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4422"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09414"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.26"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "2.102.51"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.731"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.166"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.01"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.58"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06669"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.220"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").Value = "30.078.47"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.61"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.330"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "2.349.16"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.558"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.21"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.168"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.640"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.216"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.962"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.191"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.09"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02562"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06778"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2275"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6930"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.49"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.316"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6633"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.278"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.631"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000352"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.220"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.91"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07212"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.26%  "
